$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-empty Maximum/Minimum bitrate values for rows 20-41
$ws.Range("C20").Value = 17000
$ws.Range("D20").Value = 19200

$ws.Range("C21").Value = 26000
$ws.Range("D21").Value = 31000

$ws.Range("C22").Value = 34000
$ws.Range("D22").Value = 38300

$ws.Range("C23").Value = 50000
$ws.Range("D23").Value = 59000

$ws.Range("C24").Value = 52000
$ws.Range("D24").Value = 65000

$ws.Range("C25").Value = 53000
$ws.Range("D25").Value = 58000

$ws.Range("C26").Value = 77000
$ws.Range("D26").Value = 89000

$ws.Range("C27").Value = 152000
$ws.Range("D27").Value = 175000

$ws.Range("C28").Value = 158000
$ws.Range("D28").Value = 178000

$ws.Range("C29").Value = 12000
$ws.Range("D29").Value = 14000

$ws.Range("C30").Value = 6000
$ws.Range("D30").Value = 14000

$ws.Range("C31").Value = 4000
$ws.Range("D31").Value = 11000

$ws.Range("C32").Value = 5000
$ws.Range("D32").Value = 14000

$ws.Range("C33").Value = 3000
$ws.Range("D33").Value = 14000

$ws.Range("C34").Value = 2000
$ws.Range("D34").Value = 5000

$ws.Range("C35").Value = 2000
$ws.Range("D35").Value = 5000

$ws.Range("C36").Value = 2000
$ws.Range("D36").Value = 5000

$ws.Range("C37").Value = 800
$ws.Range("D37").Value = 4000

$ws.Range("C38").Value = 800
$ws.Range("D38").Value = 2000

$ws.Range("C39").Value = 1000
$ws.Range("D39").Value = 2000

$ws.Range("C40").Value = 16000
$ws.Range("D40").Value = 19000

$ws.Range("C41").Value = 16000
$ws.Range("D41").Value = 18000

# Rows 42-54 have no numeric spec available; mark as "\"
$ws.Range("C42:D54").Value = "\"

# Update the saved selection/active cell as recorded in the sheet view
$ws.Range("D10").Select()
